$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 22:41 (only first 20 question rows remain)
$ws.Rows("22:41").Delete()

# Update the answer column (B2:B21) with the new values
$answers = @("d","a","d","a","c","b","a","c","c","a","d","c","d","a","b","c","d","b","a","c")
for ($i = 0; $i -lt $answers.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $answers[$i]
}

# Update the selected cell to B8
$ws.Range("B8").Select()
